# ---------------------------------------------------------------------------
# feat: add 2022-Q4 data
#
# 1. "总计" sheet: insert a new row (2022-Q4) at the top of the data table,
#    pushing the existing quarters (2022-Q3, 2022-Q2, 2022-Q1) down by one row.
# 2. Workbook: insert a brand new worksheet named "2022-Q4" right after
#    "总计" and before "2022-Q3", holding the per-fund breakdown for the
#    new quarter (same shape/formatting as the existing quarter sheets).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: update the "总计" (summary) sheet - shift rows 2..4 down to 3..5,
# preserving formatting, then write the new 2022-Q4 row at row 2.
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)

for ($r = 4; $r -ge 2; $r--) {
    $dst = $r + 1
    $srcRange = "A" + $r + ":D" + $r
    $dstRange = "A" + $dst + ":D" + $dst
    $wsTotal.Range($srcRange).Copy($wsTotal.Range($dstRange))
}

$wsTotal.Range("A2").Value2 = 0
$wsTotal.Range("B2").Value2 = "2022-Q4"
$wsTotal.Range("C2").Value2 = 10
$wsTotal.Range("D2").Value2 = 0.67

# ---------------------------------------------------------------------------
# Step 2: insert the new "2022-Q4" worksheet right before the current
# "2022-Q3" sheet (which sits at index 2 before insertion).
# ---------------------------------------------------------------------------
$wsQ3Before = $wb.Worksheets.Item(2)
$newSheet = $wb.Worksheets.Add($wsQ3Before)
$newSheet.Name = "2022-Q4"

# After the insert, object references shift with position, so re-fetch
# the sheets we need by their (now updated) index.
$wsQ4 = $wb.Worksheets.Item(2)
$wsQ3 = $wb.Worksheets.Item(3)

# Copy the template range (header row + its 8 data rows) from the existing
# "2022-Q3" sheet so the new sheet starts out with the exact same
# layout/styles (borders, bold header, centered index column...). The
# template only has 9 rows, but the new quarter needs 11 (header + 10
# funds), so the last template data row is duplicated twice more to cover
# rows 10 and 11, preserving their formatting as well.
$wsQ3.Range("A1:H9").Copy($wsQ4.Range("A1:H9"))
$wsQ3.Range("A9:H9").Copy($wsQ4.Range("A10:H10"))
$wsQ3.Range("A9:H9").Copy($wsQ4.Range("A11:H11"))

# Columns B..G hold text values (fund code, name, scale, position figures)
# even though some look numeric ("014269", "2.65", ...). Force them to the
# Text number format before writing so leading zeros / exact text is kept
# (matches the source data which stores these as inline strings).
$wsQ4.Range("B2:G11").NumberFormat = "@"

# ---------------------------------------------------------------------------
# Step 3: populate the fund-level data rows for 2022-Q4.
# ---------------------------------------------------------------------------
$rows = @(
    @(0, "014269", "嘉实北交所精选两年定期混合A", "2.65", "94.48", "6.47", "0.1715", 5),
    @(1, "014279", "汇添富北交所创新精选两年定开混合A", "3.06", "94.24", "5.33", "0.1631", 6),
    @(2, "014275", "易方达北交所精选两年定开混合A", "3.47", "70.57", "2.37", "0.0822", 10),
    @(3, "014663", "富国创新发展两年定期开放混合A", "2.26", "82.02", "3.61", "0.0816", 2),
    @(4, "014271", "大成北交所两年定开混合A", "3.24", "68.93", "1.97", "0.0638", 10),
    @(5, "014270", "嘉实北交所精选两年定期混合C", "0.52", "94.48", "6.47", "0.0336", 5),
    @(6, "014280", "汇添富北交所创新精选两年定开混合C", "0.48", "94.24", "5.33", "0.0256", 6),
    @(7, "014276", "易方达北交所精选两年定开混合C", "0.90", "70.57", "2.37", "0.0213", 10),
    @(8, "014272", "大成北交所两年定开混合C", "0.77", "68.93", "1.97", "0.0152", 10),
    @(9, "014664", "富国创新发展两年定期开放混合C", "0.33", "82.02", "3.61", "0.0119", 2)
)

$rowIdx = 2
foreach ($row in $rows) {
    $wsQ4.Range("A" + $rowIdx).Value2 = $row[0]
    $wsQ4.Range("B" + $rowIdx).Value2 = $row[1]
    $wsQ4.Range("C" + $rowIdx).Value2 = $row[2]
    $wsQ4.Range("D" + $rowIdx).Value2 = $row[3]
    $wsQ4.Range("E" + $rowIdx).Value2 = $row[4]
    $wsQ4.Range("F" + $rowIdx).Value2 = $row[5]
    $wsQ4.Range("G" + $rowIdx).Value2 = $row[6]
    $wsQ4.Range("H" + $rowIdx).Value2 = $row[7]
    $rowIdx = $rowIdx + 1
}

# Re-select the "2022-Q1" sheet (now shifted to index 5) so the originally
# selected tab stays the same logical sheet after the insertion - Add()
# above made the brand-new sheet the active one as a side effect.
$wsQ1 = $wb.Worksheets.Item(5)
$wsQ1.Activate()

Write-Host "2022-Q4 sheet populated."
